$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1637.75
$ws.Range("I111").Value = 1029
$ws.Range("J111").Value = 1840.6666
$ws.Range("K111").Value = 3087
$ws.Range("L111").Value = 5521.9998
$ws.Range("M111").Value = -20
$ws.Range("N111").Value = -11655.9998
$ws.Range("H137").Value = 1385.5758
$ws.Range("I137").Value = 1230.2084
$ws.Range("J137").Value = 1799.8889
$ws.Range("K137").Value = 3690.6252
$ws.Range("L137").Value = 5399.6667
$ws.Range("M137").Value = -1140.6252
$ws.Range("N137").Value = -10499.6667
$ws.Range("H138").Value = 1662.6
$ws.Range("I138").Value = 1195.7805
$ws.Range("J138").Value = 2669.9473
$ws.Range("K138").Value = 3587.3415
$ws.Range("L138").Value = 8009.841899999999
$ws.Range("M138").Value = 1552.6585
$ws.Range("N138").Value = -18289.8419
$ws.Range("H141").Value = 10478.75
$ws.Range("I141").Value = 1621.875
$ws.Range("J141").Value = 28192.5
$ws.Range("K141").Value = 4865.625
$ws.Range("L141").Value = 84577.5
$ws.Range("M141").Value = 314.375
$ws.Range("N141").Value = -94937.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2211.8965
$ws.Range("I61").Value = 2111.0908
$ws.Range("J61").Value = 2528.7144
$ws.Range("K61").Value = 2111.0908
$ws.Range("L61").Value = 2528.7144
$ws.Range("M61").Value = -1899.0908
$ws.Range("N61").Value = -2952.7144
$ws.Range("H74").Value = 1281.5
$ws.Range("I74").Value = 773.0909
$ws.Range("K74").Value = 773.0909
$ws.Range("M74").Value = 100.9091
$ws.Range("H77").Value = 1281.5
$ws.Range("I77").Value = 773.0909
$ws.Range("K77").Value = 3865.4545
$ws.Range("M77").Value = 502.5454999999997
$ws.Range("H122").Value = 3378.2856
$ws.Range("I122").Value = 3308.7273
$ws.Range("J122").Value = 3633.3333
$ws.Range("K122").Value = 9926.1819
$ws.Range("L122").Value = 10899.9999
$ws.Range("M122").Value = -7476.1819
$ws.Range("N122").Value = -15799.9999
$ws.Range("H132").Value = 4475.5347
$ws.Range("I132").Value = 4890.933
$ws.Range("J132").Value = 3516.923
$ws.Range("K132").Value = 14672.799
$ws.Range("L132").Value = 10550.769
$ws.Range("M132").Value = -12142.799
$ws.Range("N132").Value = -15610.769
$ws.Range("H136").Value = 2211.8965
$ws.Range("I136").Value = 2111.0908
$ws.Range("J136").Value = 2528.7144
$ws.Range("K136").Value = 6333.2724
$ws.Range("L136").Value = 7586.1432
$ws.Range("M136").Value = -3783.2724
$ws.Range("N136").Value = -12686.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 48811.09
$ws.Range("I107").Value = 65513.875
$ws.Range("J107").Value = 4270.3335
$ws.Range("K107").Value = 65513.875
$ws.Range("L107").Value = 4270.3335
$ws.Range("M107").Value = -63593.875
$ws.Range("N107").Value = -8110.3335
$ws.Range("H134").Value = 2462.2
$ws.Range("I134").Value = 2284.9167
$ws.Range("J134").Value = 3171.3333
$ws.Range("K134").Value = 6854.750100000001
$ws.Range("L134").Value = 9513.999899999999
$ws.Range("M134").Value = -4319.750100000001
$ws.Range("N134").Value = -14583.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3167.9048
$ws.Range("I31").Value = 2038.2142
$ws.Range("J31").Value = 5427.2856
$ws.Range("K31").Value = 2038.2142
$ws.Range("L31").Value = 5427.2856
$ws.Range("M31").Value = -1743.2142
$ws.Range("N31").Value = -6017.2856
$ws.Range("H34").Value = 3167.9048
$ws.Range("I34").Value = 2038.2142
$ws.Range("J34").Value = 5427.2856
$ws.Range("K34").Value = 2038.2142
$ws.Range("L34").Value = 5427.2856
$ws.Range("M34").Value = -1836.2142
$ws.Range("N34").Value = -5831.2856
$ws.Range("H58").Value = 1279162.8
$ws.Range("I58").Value = 1951354.1
$ws.Range("J58").Value = 1999.2
$ws.Range("K58").Value = 1951354.1
$ws.Range("L58").Value = 1999.2
$ws.Range("M58").Value = -1951151.1
$ws.Range("N58").Value = -2405.2
$ws.Range("H105").Value = 2018.4166
$ws.Range("I105").Value = 2151.25
$ws.Range("J105").Value = 1752.75
$ws.Range("K105").Value = 2151.25
$ws.Range("L105").Value = 1752.75
$ws.Range("M105").Value = -404.25
$ws.Range("N105").Value = -5246.75
$ws.Range("H132").Value = 424289.75
$ws.Range("I132").Value = 541908.4
$ws.Range("K132").Value = 1625725.2
$ws.Range("M132").Value = -1623195.2
$ws.Range("H134").Value = 1948.8948
$ws.Range("I134").Value = 1463.5927
$ws.Range("K134").Value = 4390.7781
$ws.Range("M134").Value = -1855.7781
$ws.Range("H136").Value = 1279162.8
$ws.Range("I136").Value = 1951354.1
$ws.Range("J136").Value = 1999.2
$ws.Range("K136").Value = 5854062.300000001
$ws.Range("L136").Value = 5997.6
$ws.Range("M136").Value = -5851512.300000001
$ws.Range("N136").Value = -11097.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2272.0908
$ws.Range("I5").Value = 2272.0908
$ws.Range("K5").Value = 6816.2724
$ws.Range("M5").Value = -6704.2724
$ws.Range("H12").Value = 1207786.5
$ws.Range("J12").Value = 1288305.5
$ws.Range("L12").Value = 3864916.5
$ws.Range("N12").Value = -3865262.5
$ws.Range("H135").Value = 2272.0908
$ws.Range("I135").Value = 2272.0908
$ws.Range("K135").Value = 20448.8172
$ws.Range("M135").Value = -17913.8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3292.8333
$ws.Range("I102").Value = 3222.2222
$ws.Range("J102").Value = 3504.6667
$ws.Range("K102").Value = 3222.2222
$ws.Range("L102").Value = 3504.6667
$ws.Range("M102").Value = -1600.2222
$ws.Range("N102").Value = -6748.6667
$ws.Range("H126").Value = 3063.25
$ws.Range("I126").Value = 2979.111
$ws.Range("J126").Value = 3171.4285
$ws.Range("K126").Value = 8937.332999999999
$ws.Range("L126").Value = 9514.2855
$ws.Range("M126").Value = -6467.332999999999
$ws.Range("N126").Value = -14454.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3234.4746
$ws.Range("I132").Value = 3173.9092
$ws.Range("J132").Value = 3412.1333
$ws.Range("K132").Value = 9521.7276
$ws.Range("L132").Value = 10236.3999
$ws.Range("M132").Value = -6991.7276
$ws.Range("N132").Value = -15296.3999
$ws.Range("H133").Value = 44629.152
$ws.Range("J133").Value = 44629.152
$ws.Range("L133").Value = 44629.152
$ws.Range("N133").Value = -49689.152
$ws.Range("H136").Value = 2553103
$ws.Range("I136").Value = 4033572
$ws.Range("J136").Value = 3406.111
$ws.Range("K136").Value = 12100716
$ws.Range("L136").Value = 10218.333
$ws.Range("M136").Value = -12098166
$ws.Range("N136").Value = -15318.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2715.9429
$ws.Range("I132").Value = 2287.0667
$ws.Range("J132").Value = 3037.6
$ws.Range("K132").Value = 6861.2001
$ws.Range("L132").Value = 9112.799999999999
$ws.Range("M132").Value = -4331.2001
$ws.Range("N132").Value = -14172.8
$ws.Range("H136").Value = 1728.3077
$ws.Range("I136").Value = 1377.9048
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 4133.7144
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -1583.7144
$ws.Range("N136").Value = -14700
